$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "244.87") would be
# auto-coerced by Excel into a numeric cell, which would lose formatting such
# as trailing zeros ("0.660" -> 0.66). Force these to stay text the way Excel
# does (leading apostrophe) and then reset the style to Normal so no stray
# "quote prefix" cell style lingers on the cell (matches original formatting).

$ws.Range("D2").Value = "36.433.41"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "2.006.98"
$ws.Range("E3").Value = "  +5.60%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'244.87"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").Value = "'0.660"
$ws.Range("E6").Value = "  -5.00%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'44.73"
$ws.Range("E8").Value = "  +3.52%  "
$ws.Range("D9").Value = "'61.43"
$ws.Range("E9").Value = "  +9.42%  "
$ws.Range("D10").Value = "'0.368"
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("D11").Value = "'0.0712"
$ws.Range("E11").Value = "  -6.01%  "
$ws.Range("D12").Value = "'0.0982"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "'14.56"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "2.299.43"
$ws.Range("E14").Value = "  +5.85%  "
$ws.Range("D15").Value = "'0.811"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "2.010.54"
$ws.Range("E16").Value = "  +5.40%  "
$ws.Range("D17").Value = "'4.88"
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("D18").Value = "36.387.28"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").Value = "'71.28"
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("D20").Value = "0.0₃0814"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").Value = "'12.82"
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("D22").Value = "'236.63"
$ws.Range("E22").Value = "  -3.49%  "
$ws.Range("D23").Value = "'4.87"
$ws.Range("E23").Value = "  -6.54%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  -10.41%  "
$ws.Range("D26").Value = "'165.20"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "'8.60"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "'19.60"
$ws.Range("E28").Value = "  +6.77%  "
$ws.Range("E29").Value = "  -10.83%  "
$ws.Range("E30").Value = "  -5.73%  "
$ws.Range("D31").Value = "'22.27"
$ws.Range("E31").Value = "  +62.63%  "
$ws.Range("D32").Value = "'4.37"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "'0.0584"
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -6.47%  "
$ws.Range("D37").Value = "'2.14"
$ws.Range("E37").Value = "  +8.96%  "
$ws.Range("D38").Value = "'0.0808"
$ws.Range("E38").Value = "  +8.78%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  -9.63%  "
$ws.Range("D41").Value = "'0.0215"
$ws.Range("E41").Value = "  -4.36%  "
$ws.Range("D42").Value = "'95.77"
$ws.Range("E42").Value = "  -3.70%  "
$ws.Range("D43").Value = "'1.10"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "  +15.80%  "
$ws.Range("D45").Value = "'15.96"
$ws.Range("E45").Value = "  -6.43%  "
$ws.Range("D46").Value = "1.316.36"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").Value = "'0.0816"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "2.195.20"
$ws.Range("E49").Value = "  +5.84%  "
$ws.Range("D50").Value = "'2.19"
$ws.Range("E50").Value = "  -7.77%  "
$ws.Range("D51").Value = "'3.86"
$ws.Range("E51").Value = "  +14.66%  "

# Clear the auto-applied "quote prefix" style from the forced-text cells so
# their cell formatting matches the original (unstyled) cells.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
